$d = $word.ActiveDocument

# Move to the end of the document and add a new paragraph
$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()
$end.Collapse(0)
$end.InsertAfter("Modificación o cambio realizado en el documento ")

Write-Output "done"
